# "added MP Tower striplook" - builds the Tower strip-board layout:
#   Row 2: title "Tower"
#   Row 3: ground/line-up strip (blue fill)
#   Row 4: arrival/final strip (orange fill), except the departure handover cell (blue fill)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Cell values ----------------------------------------------------------
# (entered in the original authoring order so the shared-string table lines
# up with the source workbook)
$ws.Range("A2").Value = "Tower"

$ws.Range("A4").Value = "ON FINAL"
$ws.Range("B4").Value = "ARRIVAL 07L/25R"
$ws.Range("C4").Value = "ARRIVAL 07C/25C"

$ws.Range("A3").Value = "HANDOVER FROM GROUND"
$ws.Range("E4").Value = "HANDOVER TO DEPARTURE"

$ws.Range("B3").Value = "LINE UP 07L/25R"
$ws.Range("C3").Value = "LINE UP 07C/25C"
$ws.Range("E3").Value = "LINE UP 18"
$ws.Range("D3").Value = "LINE UP 07R/25L"

$ws.Range("D4").Value = "ARRIVAL 07R/25L"

# ---- Fill colors ----------------------------------------------------------
# Blue = "Blue, Accent 1, Lighter 80%" (theme 4, tint ~0.8 -> RGB DEEBF7)
# Orange = "Orange, Accent 2, Lighter 80%" (theme 5, tint ~0.8 -> RGB FBE5D6)
$blue = 16247774
$orange = 14083579

$ws.Range("A3:E3").Interior.Color = $blue
$ws.Range("A4:D4").Interior.Color = $orange
$ws.Range("E4").Interior.Color = $blue

# ---- Column widths (best-fit to the longest label in each column) --------
$ws.Columns.Item(1).ColumnWidth = 24.59
$ws.Columns.Item(2).ColumnWidth = 17.76
$ws.Columns.Item(3).ColumnWidth = 18.09
$ws.Columns.Item(4).ColumnWidth = 18.09
$ws.Columns.Item(5).ColumnWidth = 24.09

# ---- Selection / view ------------------------------------------------------
$null = $ws.Range("D5").Select()
